$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1) Update the two raw-data-path cell values (the content change behind the
#    commit "update lake and the alum data path").
# ---------------------------------------------------------------------------

# B7: ALUM_2020-IUCNGET raw data path -> newer ABARES land-use raster
$ws.Range("B7").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2010_11_to_2020_21_prerelease3_20240809\NLUM_v7p3_ALUMV8_250m_2020_21_alb\NLUM_v7p3_ALUMV8_250m_2020_21_alb.tif"

# B5: Lacustrine-IUCNGET raw data path -> newer lakes shapefile (no longer
# overlapping with ALUM)
$ws.Range("B5").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\Lakes_NEAP_20240808_NoOverlapWithALUM.shp"

# ---------------------------------------------------------------------------
# 2) Because the raw paths behind B5 and B7 changed, the (now stale)
#    hyperlinks that used to be attached to those two cells are removed. The
#    cells keep looking like links (the "Hyperlink" cell style is restored
#    below), they are simply no longer clickable. All other hyperlinks on
#    the sheet are left exactly as they were.
#
#    This interop runtime's Hyperlinks.Delete() only works at the worksheet
#    level (it always clears every hyperlink on the sheet, even when invoked
#    through a single range/item, and individual Hyperlink.Delete() calls are
#    a no-op), so the surviving hyperlinks are recreated from scratch - in
#    their original order - right after clearing. Adding a hyperlink also
#    resets the cell's formatting, so each cell's original "Hyperlink" style
#    (and, where relevant, vertical-centered alignment) is reapplied
#    immediately afterwards.
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("I5"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\Lacustrine_EPSG3577_250m.tif") | Out-Null
$ws.Range("I5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I6"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\Estuarine_EPSG3577_250m.tif") | Out-Null
$ws.Range("I6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I2"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\NEAP_BSU_EPSG3577_250m.tif") | Out-Null
$ws.Range("I2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I3"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\Marine_EPSG3577_250m.tif") | Out-Null
$ws.Range("I3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B2"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\BSU\outputs\BSU_NEAP\BSU_NEAP_epsg3577_250m.tif") | Out-Null
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B2").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("J5"), "https://github.com/CSIRO-enviro-informatics/ecosystem-typology/raw/main/crosswalks/Geofabric-IUCNGET/Lacustrine-IUCNGET.xlsx") | Out-Null
$ws.Range("J5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("I4"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\Terrestrial_Extant_EPSG3577_250m.tif") | Out-Null
$ws.Range("I4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B4"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\NVIS_IUCNGET_DK_20240801.tif") | Out-Null
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B4").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("I7"), "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\ALUM_2020_EPSG3577_250m.tif") | Out-Null
$ws.Range("I7").Style = "Hyperlink"
